$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "leak" edge list gained a batch of reverse-direction edges (f<a>t<b>Leak
# alongside the already-present f<b>t<a>Leak), so the header row's label set
# grew from 55 to 82 unique strings. Because the new labels were inserted in
# the middle of the logical ordering (not merely appended), the label shown
# in every header cell from column Z onward shifts to the next one in that
# ordering, and the row grows out to column CE to fit the newcomers.
$headers = @(
    "N0", "N1", "N2", "PWFn2", "N3", "N4", "PWFn4", "N5",
    "N6", "N7", "N8", "N9", "N11", "N12", "N13", "N14",
    "N15", "N16", "PWFn16", "N17", "N18", "N19", "N20", "PWFn20",
    "f0t1Leak", "f0t7Leak", "f1t0Leak", "f1t8Leak", "f1t2Leak", "f2t1Leak", "f2t9Leak", "f2t3Leak",
    "f3t2Leak", "f3t4Leak", "f4t3Leak", "f4t11Leak", "f4t5Leak", "f5t4Leak", "f5t6Leak", "f5t12Leak",
    "f6t5Leak", "f6t13Leak", "f7t0Leak", "f7t8Leak", "f7t14Leak", "f8t7Leak", "f8t1Leak", "f8t9Leak",
    "f9t2Leak", "f9t8Leak", "f9t16Leak", "f11t4Leak", "f11t18Leak", "f11t12Leak", "f12t5Leak", "f12t11Leak",
    "f12t13Leak", "f12t19Leak", "f13t6Leak", "f13t12Leak", "f13t20Leak", "f14t7Leak", "f14t15Leak", "f15t14Leak",
    "f15t16Leak", "f16t9Leak", "f16t15Leak", "f16t17Leak", "f17t16Leak", "f17t18Leak", "f18t17Leak", "f18t11Leak",
    "f18t19Leak", "f19t18Leak", "f19t12Leak", "f19t20Leak", "f20t13Leak", "f20t19Leak", "RES2", "RES4",
    "RES16", "RES20"
)

$firstCol = 2          # B
$oldLastCol = 56        # BD (old last header column)
$lastCol = $firstCol + $headers.Length - 1   # CE (57 + 82 - 1 = 83)

# Stamp bold/centered/bordered style (copied from the existing B1 header
# cell) onto the freshly added header cells beyond the old BD1 end, then
# write every header cell's text (covers both the untouched-looking ones and
# the ones whose label shifted).
for ($col = $oldLastCol + 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item(1, 2).Copy($ws.Cells.Item(1, $col))
}

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $firstCol + $i).Value = $headers[$i]
}

# Row 2: every data column (B..CE) is now exactly 0 -- the previously tiny
# floating-point residuals are cleared, and the new columns start at 0 too.
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $ws.Cells.Item(2, $col).Value = 0
}
